$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}

Replace-Text "2025-03-06 Thursday" "2025-03-07 Friday"

Replace-Text "289÷6=48, 1" "881÷5=176, 1"
Replace-Text "934÷2=467, 0" "559÷9=62, 1"
Replace-Text "173÷8=21, 5" "592÷9=65, 7"
Replace-Text "372÷9=41, 3" "822÷6=137, 0"
Replace-Text "415÷4=103, 3" "633÷2=316, 1"

Replace-Text "498÷5=99, 3" "993÷9=110, 3"
Replace-Text "986÷3=328, 2" "961÷8=120, 1"
Replace-Text "339÷6=56, 3" "908÷4=227, 0"
Replace-Text "382÷8=47, 6" "249÷6=41, 3"
Replace-Text "973÷8=121, 5" "811÷9=90, 1"

Replace-Text "847÷3=282, 1" "873÷4=218, 1"
Replace-Text "266÷9=29, 5" "686÷4=171, 2"
Replace-Text "603÷6=100, 3" "560÷6=93, 2"
Replace-Text "203÷4=50, 3" "996÷7=142, 2"
Replace-Text "391÷7=55, 6" "751÷7=107, 2"

Replace-Text "286÷4=71, 2" "622÷4=155, 2"
Replace-Text "456÷6=76, 0" "316÷3=105, 1"
Replace-Text "836÷5=167, 1" "640÷8=80, 0"
Replace-Text "942÷5=188, 2" "575÷2=287, 1"
Replace-Text "312÷8=39, 0" "445÷9=49, 4"

Replace-Text "201÷4=50, 1" "898÷3=299, 1"
Replace-Text "564÷4=141, 0" "888÷9=98, 6"
Replace-Text "509÷2=254, 1" "698÷7=99, 5"
Replace-Text "972÷8=121, 4" "623÷6=103, 5"
Replace-Text "800÷4=200, 0" "486÷5=97, 1"

Write-Output "Done applying replacements"
